$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 and F5 held a bank account number that Excel had mangled into a
# floating point number (loss of precision / scientific notation). Store
# the correct value as text instead.
$ws.Range("F2").Value = "21029112-11111112-25163151"
$ws.Range("F5").Value = "21029112-11111112-25163151"

# The selection was left on F5 when the file was saved.
$ws.Range("F5").Select()
